$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Cells.Item(10,4).Value = 44428
$ws.Cells.Item(10,11).Value = 1500
$ws.Cells.Item(10,12).Value = 1800
$ws.Cells.Item(10,13).Value = 1650
$ws.Cells.Item(10,16).Value = 550

# Row 11
$ws.Cells.Item(11,4).Value = 44342
$ws.Cells.Item(11,12).Value = 2200
$ws.Cells.Item(11,13).Value = 2100
$ws.Cells.Item(11,16).Value = 700

# Row 12
$ws.Cells.Item(12,4).Value = 44385
$ws.Cells.Item(12,11).Value = 2000
$ws.Cells.Item(12,12).Value = 2300
$ws.Cells.Item(12,13).Value = 2150
$ws.Cells.Item(12,16).Value = 717

# Row 13
$ws.Cells.Item(13,4).Value = 44295
$ws.Cells.Item(13,10).Value = 200
$ws.Cells.Item(13,11).Value = 1500
$ws.Cells.Item(13,13).Value = 1650
$ws.Cells.Item(13,16).Value = 550

# Row 14
$ws.Cells.Item(14,4).Value = 44300
$ws.Cells.Item(14,10).Value = 250
$ws.Cells.Item(14,11).Value = 1600
$ws.Cells.Item(14,13).Value = 1700
$ws.Cells.Item(14,16).Value = 567

# Row 15
$ws.Cells.Item(15,9).Value = 'Primera'
$ws.Cells.Item(15,11).Value = 1700
$ws.Cells.Item(15,12).Value = 1800
$ws.Cells.Item(15,13).Value = 1750
$ws.Cells.Item(15,16).Value = 583

# Row 16
$ws.Cells.Item(16,4).Value = 44364
$ws.Cells.Item(16,9).Value = 'Segunda'

# Row 17
$ws.Cells.Item(17,4).Value = 44435
$ws.Cells.Item(17,9).Value = 'Primera'
$ws.Cells.Item(17,10).Value = 450
$ws.Cells.Item(17,12).Value = 1300
$ws.Cells.Item(17,13).Value = 1194
$ws.Cells.Item(17,16).Value = 398

# Row 18
$ws.Cells.Item(18,4).Value = 44435
$ws.Cells.Item(18,9).Value = 'Segunda'
$ws.Cells.Item(18,11).Value = 950
$ws.Cells.Item(18,12).Value = 1000
$ws.Cells.Item(18,13).Value = 975
$ws.Cells.Item(18,16).Value = 325

# Row 19
$ws.Cells.Item(19,4).Value = 44431
$ws.Cells.Item(19,11).Value = 1000
$ws.Cells.Item(19,12).Value = 1300
$ws.Cells.Item(19,13).Value = 1150
$ws.Cells.Item(19,16).Value = 383

# Row 20
$ws.Cells.Item(20,4).Value = 44224
$ws.Cells.Item(20,10).Value = 200
$ws.Cells.Item(20,11).Value = 1400
$ws.Cells.Item(20,12).Value = 1500
$ws.Cells.Item(20,13).Value = 1450
$ws.Cells.Item(20,16).Value = 483

# Row 21
$ws.Cells.Item(21,4).Value = 44224
$ws.Cells.Item(21,10).Value = 160
$ws.Cells.Item(21,11).Value = 1000
$ws.Cells.Item(21,12).Value = 1200
$ws.Cells.Item(21,13).Value = 1100
$ws.Cells.Item(21,16).Value = 367

# Row 22
$ws.Cells.Item(22,4).Value = 44327
$ws.Cells.Item(22,10).Value = 200
$ws.Cells.Item(22,11).Value = 1400
$ws.Cells.Item(22,12).Value = 1500
$ws.Cells.Item(22,13).Value = 1450
$ws.Cells.Item(22,16).Value = 483

# Row 23
$ws.Cells.Item(23,4).Value = 44391
$ws.Cells.Item(23,9).Value = 'Primera'
$ws.Cells.Item(23,10).Value = 250
$ws.Cells.Item(23,11).Value = 1800
$ws.Cells.Item(23,12).Value = 2000
$ws.Cells.Item(23,13).Value = 1900
$ws.Cells.Item(23,16).Value = 633

# Row 24
$ws.Cells.Item(24,4).Value = 44278
$ws.Cells.Item(24,10).Value = 140
$ws.Cells.Item(24,11).Value = 2000
$ws.Cells.Item(24,12).Value = 2500
$ws.Cells.Item(24,13).Value = 2250
$ws.Cells.Item(24,16).Value = 750

# Row 25
$ws.Cells.Item(25,4).Value = 44278
$ws.Cells.Item(25,10).Value = 200
$ws.Cells.Item(25,11).Value = 1500
$ws.Cells.Item(25,12).Value = 1800
$ws.Cells.Item(25,13).Value = 1650
$ws.Cells.Item(25,16).Value = 550

# Row 26
$ws.Cells.Item(26,4).Value = 44417
$ws.Cells.Item(26,10).Value = 250
$ws.Cells.Item(26,11).Value = 1800
$ws.Cells.Item(26,12).Value = 2000
$ws.Cells.Item(26,13).Value = 1900
$ws.Cells.Item(26,16).Value = 633

# Row 27
$ws.Cells.Item(27,4).Value = 44417
$ws.Cells.Item(27,9).Value = 'Segunda'
$ws.Cells.Item(27,10).Value = 200
$ws.Cells.Item(27,11).Value = 1500
$ws.Cells.Item(27,12).Value = 1600
$ws.Cells.Item(27,13).Value = 1550
$ws.Cells.Item(27,14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(27,16).Value = 517
$ws.Cells.Item(27,17).Value = 3

# Row 28
$ws.Cells.Item(28,4).Value = 44161
$ws.Cells.Item(28,10).Value = 200
$ws.Cells.Item(28,11).Value = 600
$ws.Cells.Item(28,12).Value = 700
$ws.Cells.Item(28,13).Value = 650
$ws.Cells.Item(28,16).Value = 217

# Row 29
$ws.Cells.Item(29,4).Value = 44161
$ws.Cells.Item(29,9).Value = 'Segunda'
$ws.Cells.Item(29,10).Value = 250
$ws.Cells.Item(29,11).Value = 500
$ws.Cells.Item(29,12).Value = 600
$ws.Cells.Item(29,13).Value = 550
$ws.Cells.Item(29,16).Value = 183

# Row 30
$ws.Cells.Item(30,4).Value = 44333
$ws.Cells.Item(30,9).Value = 'Primera'
$ws.Cells.Item(30,10).Value = 200
$ws.Cells.Item(30,11).Value = 1500
$ws.Cells.Item(30,12).Value = 1700
$ws.Cells.Item(30,13).Value = 1600
$ws.Cells.Item(30,16).Value = 533

# Row 31
$ws.Cells.Item(31,4).Value = 44398
$ws.Cells.Item(31,10).Value = 300
$ws.Cells.Item(31,11).Value = 1700
$ws.Cells.Item(31,12).Value = 1800
$ws.Cells.Item(31,13).Value = 1750
$ws.Cells.Item(31,16).Value = 583

# Row 32
$ws.Cells.Item(32,1).Value = 1
$ws.Cells.Item(32,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(32,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(32,4).Value = 44270
$ws.Cells.Item(32,5).Value = 15
$ws.Cells.Item(32,6).Value = 100112009
$ws.Cells.Item(32,7).Value = 'Acelga'
$ws.Cells.Item(32,8).Value = 'Sin especificar'
$ws.Cells.Item(32,9).Value = 'Primera'
$ws.Cells.Item(32,10).Value = 100
$ws.Cells.Item(32,11).Value = 1800
$ws.Cells.Item(32,12).Value = 2000
$ws.Cells.Item(32,13).Value = 1900
$ws.Cells.Item(32,14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(32,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(32,16).Value = 633
$ws.Cells.Item(32,17).Value = 3
$ws.Cells.Item(32,18).Value = 'Hortaliza'
$ws.Cells.Item(32,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 33
$ws.Cells.Item(33,1).Value = 1
$ws.Cells.Item(33,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(33,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(33,4).Value = 44270
$ws.Cells.Item(33,5).Value = 15
$ws.Cells.Item(33,6).Value = 100112009
$ws.Cells.Item(33,7).Value = 'Acelga'
$ws.Cells.Item(33,8).Value = 'Sin especificar'
$ws.Cells.Item(33,9).Value = 'Segunda'
$ws.Cells.Item(33,10).Value = 100
$ws.Cells.Item(33,11).Value = 1200
$ws.Cells.Item(33,12).Value = 1500
$ws.Cells.Item(33,13).Value = 1350
$ws.Cells.Item(33,14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(33,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(33,16).Value = 450
$ws.Cells.Item(33,17).Value = 3
$ws.Cells.Item(33,18).Value = 'Hortaliza'
$ws.Cells.Item(33,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 34
$ws.Cells.Item(34,1).Value = 1
$ws.Cells.Item(34,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(34,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(34,4).Value = 44376
$ws.Cells.Item(34,5).Value = 15
$ws.Cells.Item(34,6).Value = 100112009
$ws.Cells.Item(34,7).Value = 'Acelga'
$ws.Cells.Item(34,8).Value = 'Sin especificar'
$ws.Cells.Item(34,9).Value = 'Primera'
$ws.Cells.Item(34,10).Value = 340
$ws.Cells.Item(34,11).Value = 1400
$ws.Cells.Item(34,12).Value = 1500
$ws.Cells.Item(34,13).Value = 1471
$ws.Cells.Item(34,14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(34,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(34,16).Value = 490
$ws.Cells.Item(34,17).Value = 3
$ws.Cells.Item(34,18).Value = 'Hortaliza'
$ws.Cells.Item(34,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 35
$ws.Cells.Item(35,1).Value = 1
$ws.Cells.Item(35,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(35,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(35,4).Value = 44432
$ws.Cells.Item(35,5).Value = 15
$ws.Cells.Item(35,6).Value = 100112009
$ws.Cells.Item(35,7).Value = 'Acelga'
$ws.Cells.Item(35,8).Value = 'Sin especificar'
$ws.Cells.Item(35,9).Value = 'Primera'
$ws.Cells.Item(35,10).Value = 200
$ws.Cells.Item(35,11).Value = 1200
$ws.Cells.Item(35,12).Value = 1300
$ws.Cells.Item(35,13).Value = 1250
$ws.Cells.Item(35,14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(35,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(35,16).Value = 417
$ws.Cells.Item(35,17).Value = 3
$ws.Cells.Item(35,18).Value = 'Hortaliza'
$ws.Cells.Item(35,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 36
$ws.Cells.Item(36,1).Value = 1
$ws.Cells.Item(36,2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(36,3).Value = 'Arica y Parinacota'
$ws.Cells.Item(36,4).Value = 44432
$ws.Cells.Item(36,5).Value = 15
$ws.Cells.Item(36,6).Value = 100112009
$ws.Cells.Item(36,7).Value = 'Acelga'
$ws.Cells.Item(36,8).Value = 'Sin especificar'
$ws.Cells.Item(36,9).Value = 'Segunda'
$ws.Cells.Item(36,10).Value = 200
$ws.Cells.Item(36,11).Value = 950
$ws.Cells.Item(36,12).Value = 1000
$ws.Cells.Item(36,13).Value = 975
$ws.Cells.Item(36,14).Value = '$/atado 2,5 a 3 kilos'
$ws.Cells.Item(36,15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(36,16).Value = 325
$ws.Cells.Item(36,17).Value = 3
$ws.Cells.Item(36,18).Value = 'Hortaliza'
$ws.Cells.Item(36,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
